$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("METALS")

$pivotSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$pivotSheet.Name = "PIVOT"

$sourceRange = $ws.Range("A1:C213")
$pc = $wb.PivotCaches().Create(1, $sourceRange)
$pt = $pc.CreatePivotTable($pivotSheet.Range("B4"), "MetalsPivot", $false, $false)

$avgField = $pt.PivotFields("avg_price")
Write-Host "avg orientation before:" $avgField.Orientation
$avgField.Orientation = 0
Write-Host "avg orientation after:" $avgField.Orientation
Write-Host "dfs count after:" $pt.DataFields().Count
